# 14.01.21 Filled all databases excluding modifications and intact
# modifications, debugging.
# Refreshes the "analysis" sheet's run timestamp and the occupancy
# values (columns C = 5'-O, D = 3'-O) computed for sequence positions
# 1..19 (rows 9..27), matching the re-run of the occupancy analysis.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("analysis")

# Update the timestamp label in A1
$ws.Range("A1").Value = "11/01/2021 17:56"

# Row 3-4 (top summary numbers)
$ws.Range("C3").Value = 0.4586934551182387
$ws.Range("C4").Value = 0.5413065448817613

# Row 9 (seq 1): clear C9 entirely (cell becomes empty, not just zero)
$ws.Range("C9").Clear()

# Row 14 (seq 6)
$ws.Range("D14").Value = 0.9426793405802766

# Row 15 (seq 7)
$ws.Range("C15").Value = 0.2948192554489688

# Row 16 (seq 8)
$ws.Range("C16").Value = 0.3089824189933628
$ws.Range("D16").Value = 0.5904833814888236

# Row 17 (seq 9)
$ws.Range("C17").Value = 0.3091405136994379
$ws.Range("D17").Value = 0.5877238379963322

# Row 18 (seq 10)
$ws.Range("C18").Value = 0.3166188715203264
$ws.Range("D18").Value = 0.6950272937821849

# Row 19 (seq 11): C19 newly added as 0 (percentage-formatted, matching column C style)
$ws.Range("C19").Value = 0
$ws.Range("C19").NumberFormat = "0.0%"
$ws.Range("D19").Value = 0.4799780978171799

# Row 20 (seq 12)
$ws.Range("C20").Value = 0.3676918970110236
$ws.Range("D20").Value = 0.4964041833034137

# Row 21 (seq 13)
$ws.Range("C21").Value = 0.5694399995102261
$ws.Range("D21").Value = 0.4511452428471781

# Row 22 (seq 14)
$ws.Range("C22").Value = 0.6968909473354162
$ws.Range("D22").Value = 0.1808410096349255

# Row 23 (seq 15)
$ws.Range("C23").Value = 0.7322849197478155
$ws.Range("D23").Value = 0.1557090277338459

# Row 24 (seq 16)
$ws.Range("C24").Value = 0.6805195656183084
$ws.Range("D24").Value = 0.2110274212379415

# Row 25 (seq 17)
$ws.Range("C25").Value = 0.871710627973019
$ws.Range("D25").Value = 0.2286928381609282

# Row 26 (seq 18)
$ws.Range("C26").Value = 0.9542488324553626
$ws.Range("D26").Value = 0.044303539091117

# Row 27 (seq 19)
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 0.0865902006984151

$wb.Save()
